# Weekly update: add a new price observation for Berenjena (Terminal La
# Palmera de La Serena) dated 2022-01-07, inserted above the existing
# historical rows (which all shift down by one row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 106, pushing the old rows 106:111 down to
# 107:112 (Excel's normal "insert row above" behaviour).
$ws.Rows("106:106").Insert()

# Populate the newly inserted row with the new weekly data point.
$ws.Range("A106").Value = 8
$ws.Range("B106").Value = "Terminal La Palmera de La Serena"
$ws.Range("C106").Value = "Coquimbo"
$ws.Range("D106").Value = 44568
$ws.Range("E106").Value = 4
$ws.Range("F106").Value = 100112001
$ws.Range("G106").Value = "Berenjena"
$ws.Range("H106").Value = "Sin especificar"
$ws.Range("I106").Value = "Primera"
$ws.Range("J106").Value = 700
$ws.Range("K106").Value = 8000
$ws.Range("L106").Value = 9000
$ws.Range("M106").Value = 8500
$ws.Range("N106").Value = "`$/caja 60 unidades"
$ws.Range("O106").Value = "Región de Arica y Parinacota"
$ws.Range("P106").Value = 142
$ws.Range("Q106").Value = 60
$ws.Range("R106").Value = "Hortaliza"
